$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so values are not
# reinterpreted as numbers/percentages by Excel's type inference.
$ws.Range("D2:E28").NumberFormat = "@"
$ws.Range("D40:E47").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "286.65"
$ws.Cells.Item(2, 5).Value = "1.39%"
$ws.Cells.Item(3, 4).Value = "29.51"
$ws.Cells.Item(3, 5).Value = "4.06%"
$ws.Cells.Item(4, 4).Value = "5.082"
$ws.Cells.Item(4, 5).Value = "1.33%"
$ws.Cells.Item(5, 4).Value = "0.06716"
$ws.Cells.Item(5, 5).Value = "3.25%"
$ws.Cells.Item(6, 4).Value = "7.342"
$ws.Cells.Item(6, 5).Value = "1.64%"
$ws.Cells.Item(7, 4).Value = "3.441"
$ws.Cells.Item(7, 5).Value = "2.40%"
$ws.Cells.Item(8, 4).Value = "1.385"
$ws.Cells.Item(8, 5).Value = "-0.67%"
$ws.Cells.Item(9, 4).Value = "0.9201"
$ws.Cells.Item(9, 5).Value = "0.23%"
$ws.Cells.Item(10, 4).Value = "0.1597"
$ws.Cells.Item(10, 5).Value = "3.85%"
$ws.Cells.Item(11, 4).Value = "0.06783"
$ws.Cells.Item(11, 5).Value = "7.69%"
$ws.Cells.Item(12, 4).Value = "0.07702"
$ws.Cells.Item(12, 5).Value = "1.83%"
$ws.Cells.Item(13, 4).Value = "0.02929"
$ws.Cells.Item(13, 5).Value = "3.64%"
$ws.Cells.Item(14, 4).Value = "0.08989"
$ws.Cells.Item(14, 5).Value = "0.34%"
$ws.Cells.Item(15, 4).Value = "0.001591"
$ws.Cells.Item(15, 5).Value = "0.08%"
$ws.Cells.Item(16, 4).Value = "0.04461"
$ws.Cells.Item(16, 5).Value = "0.70%"
$ws.Cells.Item(17, 4).Value = "0.0006444"
$ws.Cells.Item(17, 5).Value = "1.51%"
$ws.Cells.Item(18, 4).Value = "0.006278"
$ws.Cells.Item(18, 5).Value = "2.91%"
$ws.Cells.Item(19, 4).Value = "3.455"
$ws.Cells.Item(19, 5).Value = "0.29%"
$ws.Cells.Item(20, 4).Value = "2.227"
$ws.Cells.Item(20, 5).Value = "-0.63%"
$ws.Cells.Item(21, 4).Value = "0.3198"
$ws.Cells.Item(21, 5).Value = "0.50%"
$ws.Cells.Item(22, 4).Value = "0.1310"
$ws.Cells.Item(22, 5).Value = "-0.92%"
$ws.Cells.Item(23, 4).Value = "4.070"
$ws.Cells.Item(23, 5).Value = "2.47%"
$ws.Cells.Item(24, 5).Value = "2.43%"
$ws.Cells.Item(25, 4).Value = "0.001195"
$ws.Cells.Item(25, 5).Value = "1.22%"
$ws.Cells.Item(26, 4).Value = "0.004124"
$ws.Cells.Item(26, 5).Value = "-7.44%"
$ws.Cells.Item(27, 4).Value = "0.0001199"
$ws.Cells.Item(27, 5).Value = "0.00%"
$ws.Cells.Item(28, 5).Value = "-0.04%"
$ws.Cells.Item(40, 4).Value = "0.04286"
$ws.Cells.Item(40, 5).Value = "4.38%"
$ws.Cells.Item(41, 4).Value = "0.006742"
$ws.Cells.Item(41, 5).Value = "1.26%"
$ws.Cells.Item(42, 4).Value = "0.1240"
$ws.Cells.Item(42, 5).Value = "0.53%"
$ws.Cells.Item(43, 4).Value = "0.002237"
$ws.Cells.Item(43, 5).Value = "5.17%"
$ws.Cells.Item(44, 4).Value = "0.01196"
$ws.Cells.Item(44, 5).Value = "4.24%"
$ws.Cells.Item(45, 4).Value = "0.00005684"
$ws.Cells.Item(45, 5).Value = "1.04%"
$ws.Cells.Item(46, 4).Value = "1.968"
$ws.Cells.Item(46, 5).Value = "0.89%"
$ws.Cells.Item(47, 4).Value = "0.01505"
$ws.Cells.Item(47, 5).Value = "-18.66%"
